$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114 (shifts existing rows 114-148 down to 115-149)
$ws.Rows(114).Insert()

# Populate the new row 114 with this week's new data entry
$ws.Range("A114").Value = 6
$ws.Range("B114").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C114").Value = "Metropolitana"
$ws.Range("D114").Value = 44588
$ws.Range("D114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E114").Value = 13
$ws.Range("F114").Value = 100112029
$ws.Range("G114").Value = "Orégano"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 34
$ws.Range("K114").Value = 8000
$ws.Range("L114").Value = 9000
$ws.Range("M114").Value = 8441
$ws.Range("N114").Value = "$/docena de atados"
$ws.Range("O114").Value = "Región Metropolitana"
$ws.Range("P114").Value = 2814
$ws.Range("Q114").Value = 3
$ws.Range("R114").Value = "Hortaliza"
